$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q (17) into new column R (18): header + all data rows
$lastRow = 18
$srcCol = 17  # Q
$dstCol = 18  # R

# Header (row 1) - new date label "28-jun"
$ws.Cells.Item(1, $dstCol).Value = "28-jun"

# Data rows 2..18 - mirror the numeric values from column Q
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dstCol).Value = $ws.Cells.Item($r, $srcCol).Value2
}

# Update the active selection / view to match the authored file
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("S5").Select()
